# Auto-generated Excel COM-interop edit script
# Applies the "Update automàtic: dades i banners [2026-02-09 20:20]" changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column H hold percentage text (e.g. "58%"). A direct .Value assignment
# would make Excel auto-convert these strings into numeric percentages, which
# changes both the stored type and value. Forcing a Text number format first
# keeps them as literal strings, matching the source data.
$percentCells = @("H7", "H9", "H25", "H28", "H32", "H42", "H46")
foreach ($cell in $percentCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-09 20:18:25"
$ws.Range("I2").Value = "0.8 mm"
$ws.Range("E3").Value = "2026-02-09 20:18:27"
$ws.Range("I3").Value = "1.8 mm"
$ws.Range("E4").Value = "2026-02-09 20:18:30"
$ws.Range("E5").Value = "2026-02-09 20:18:33"
$ws.Range("E6").Value = "2026-02-09 20:18:36"
$ws.Range("E7").Value = "2026-02-09 20:18:38"
$ws.Range("H7").Value = "58%"
$ws.Range("E8").Value = "2026-02-09 20:18:41"
$ws.Range("L8").Value = "36.0 km/h - 252º 19:56 TU"
$ws.Range("O8").Value = "9.1 °C"
$ws.Range("E9").Value = "2026-02-09 20:18:44"
$ws.Range("H9").Value = "81%"
$ws.Range("O9").Value = "8.6 °C"
$ws.Range("E10").Value = "2026-02-09 20:18:46"
$ws.Range("E11").Value = "2026-02-09 20:18:49"
$ws.Range("E12").Value = "2026-02-09 20:18:51"
$ws.Range("E13").Value = "2026-02-09 20:18:54"
$ws.Range("J13").Value = "1008.1 hPa"
$ws.Range("O13").Value = "3.1 °C"
$ws.Range("E14").Value = "2026-02-09 20:18:57"
$ws.Range("E15").Value = "2026-02-09 20:18:59"
$ws.Range("E16").Value = "2026-02-09 20:19:01"
$ws.Range("O16").Value = "-3.4 °C"
$ws.Range("E17").Value = "2026-02-09 20:19:04"
$ws.Range("O17").Value = "1.1 °C"
$ws.Range("E18").Value = "2026-02-09 20:19:07"
$ws.Range("E19").Value = "2026-02-09 20:19:10"
$ws.Range("E20").Value = "2026-02-09 20:19:12"
$ws.Range("E21").Value = "2026-02-09 20:19:15"
$ws.Range("E22").Value = "2026-02-09 20:19:17"
$ws.Range("E23").Value = "2026-02-09 20:19:20"
$ws.Range("I23").Value = "0.4 mm"
$ws.Range("E24").Value = "2026-02-09 20:19:22"
$ws.Range("I24").Value = "0.9 mm"
$ws.Range("J24").Value = "1008.6 hPa"
$ws.Range("O24").Value = "8.0 °C"
$ws.Range("E25").Value = "2026-02-09 20:19:24"
$ws.Range("H25").Value = "74%"
$ws.Range("E26").Value = "2026-02-09 20:19:27"
$ws.Range("O26").Value = "2.5 °C"
$ws.Range("E27").Value = "2026-02-09 20:19:30"
$ws.Range("E28").Value = "2026-02-09 20:19:32"
$ws.Range("H28").Value = "80%"
$ws.Range("E29").Value = "2026-02-09 20:19:35"
$ws.Range("O29").Value = "8.6 °C"
$ws.Range("E30").Value = "2026-02-09 20:19:37"
$ws.Range("O30").Value = "8.5 °C"
$ws.Range("E31").Value = "2026-02-09 20:19:40"
$ws.Range("E32").Value = "2026-02-09 20:19:43"
$ws.Range("H32").Value = "80%"
$ws.Range("I32").Value = "0.1 mm"
$ws.Range("O32").Value = "5.2 °C"
$ws.Range("E33").Value = "2026-02-09 20:19:45"
$ws.Range("E34").Value = "2026-02-09 20:19:48"
$ws.Range("E35").Value = "2026-02-09 20:19:51"
$ws.Range("I35").Value = "1.5 mm"
$ws.Range("E36").Value = "2026-02-09 20:19:53"
$ws.Range("E37").Value = "2026-02-09 20:19:56"
$ws.Range("E38").Value = "2026-02-09 20:19:59"
$ws.Range("E39").Value = "2026-02-09 20:20:01"
$ws.Range("E40").Value = "2026-02-09 20:20:04"
$ws.Range("E41").Value = "2026-02-09 20:20:06"
$ws.Range("E42").Value = "2026-02-09 20:20:09"
$ws.Range("H42").Value = "85%"
$ws.Range("E43").Value = "2026-02-09 20:20:12"
$ws.Range("K43").Value = "11.3 MJ/m2"
$ws.Range("E44").Value = "2026-02-09 20:20:14"
$ws.Range("O44").Value = "-3.9 °C"
$ws.Range("E45").Value = "2026-02-09 20:20:17"
$ws.Range("J45").Value = "1007.1 hPa"
$ws.Range("O45").Value = "4.0 °C"
$ws.Range("E46").Value = "2026-02-09 20:20:20"
$ws.Range("H46").Value = "73%"
$ws.Range("I46").Value = "0.7 mm"
